# Rename the three logo pictures (two Pearson logos in the footers, one
# BTec logo in the header) so that the "file name" recorded on each
# picture's document properties matches the new naming scheme:
#   - Pearson Edexcel logo (footers):  image1.png -> image2.png
#   - BTec logo (header):              image2.jpg -> image1.jpg
#
# The picture's description text (AlternativeText) is left untouched -
# only the Name (the wp:docPr / pic:cNvPr "name" attribute) changes.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer 1 : Pearson Edexcel logo -> image2.png -----------------------
$ftr1 = $sec.Footers.Item(1)
$pearsonShape1 = $ftr1.Range.InlineShapes.Item(1)
$pearsonShape1.Name = "image2.png"

# --- Footer 2 : Pearson Edexcel logo -> image2.png -----------------------
$ftr2 = $sec.Footers.Item(2)
$pearsonShape2 = $ftr2.Range.InlineShapes.Item(1)
$pearsonShape2.Name = "image2.png"

# --- Header 2 (first-page header) : BTec logo -> image1.jpg --------------
$hdr2 = $sec.Headers.Item(2)
$btecShape = $hdr2.Range.InlineShapes.Item(1)
$btecShape.Name = "image1.jpg"

Write-Host "Footer1 picture name:" $pearsonShape1.Name
Write-Host "Footer2 picture name:" $pearsonShape2.Name
Write-Host "Header2 picture name:" $btecShape.Name
